# Refresh the cryptocurrency price/volume snapshot values (GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") values are plain text in this sheet, even
# when they look numeric, so a leading apostrophe is used where needed to stop Excel from
# re-interpreting a text price (e.g. "309.68") as a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.388.72'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.374.12'
$ws.Range('E3').Value = '  +3.16%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''309.68'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '''105.35'
$ws.Range('E6').Value = '  +4.81%  '
$ws.Range('D7').Value = '''0.520'
$ws.Range('E7').Value = '  -2.89%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.518'
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('D10').Value = '''36.16'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').Value = '''53.30'
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('D13').Value = '''0.113'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').Value = '''7.00'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').Value = '2.743.50'
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').Value = '''15.61'
$ws.Range('E16').Value = '  +4.39%  '
$ws.Range('D17').Value = '2.377.00'
$ws.Range('E17').Value = '  +3.28%  '
$ws.Range('D18').Value = '''0.816'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').Value = '43.350.90'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').Value = '  -3.83%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.29'
$ws.Range('E21').Value = '  +3.49%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.0₃0920'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = '''68.32'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '''241.76'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').Value = '''2.05'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('D26').Value = '''2.62'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '''25.85'
$ws.Range('E28').Value = '  +6.13%  '
$ws.Range('D29').Value = '''3.83'
$ws.Range('E29').Value = '  -4.43%  '
$ws.Range('D30').Value = '''36.95'
$ws.Range('E30').Value = '  -4.53%  '
$ws.Range('D31').Value = '''9.59'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').Value = '''162.06'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('D34').Value = '''5.27'
$ws.Range('E34').Value = '  -1.33%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '''18.26'
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('D37').Value = '''4.81'
$ws.Range('E37').Value = '  +13.49%  '
$ws.Range('D38').Value = '''2.55'
$ws.Range('E38').Value = '  +6.57%  '
$ws.Range('D39').Value = '''3.12'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  +6.00%  '
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('D44').Value = '''2.49'
$ws.Range('E44').Value = '  +9.48%  '
$ws.Range('D45').Value = '''19.97'
$ws.Range('E45').Value = '  +3.62%  '
$ws.Range('D46').Value = '2.007.47'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '''3.15'
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('D49').Value = '''10.55'
$ws.Range('E49').Value = '  +7.36%  '
$ws.Range('D50').Value = '''58.10'
$ws.Range('E50').Value = '  +4.89%  '
$ws.Range('E51').Value = '  +1.19%  '
